$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "333.96"
Set-TextValue $ws "E2" "1.38%"
Set-TextValue $ws "G2" "12"

Set-TextValue $ws "D3" "43.89"
Set-TextValue $ws "E3" "6.13%"
Set-TextValue $ws "G3" "12"

Set-TextValue $ws "D4" "5.745"
Set-TextValue $ws "E4" "1.25%"
Set-TextValue $ws "G4" "12"

Set-TextValue $ws "D5" "0.08357"
Set-TextValue $ws "E5" "0.62%"
Set-TextValue $ws "G5" "12"

Set-TextValue $ws "D6" "8.848"
Set-TextValue $ws "E6" "0.99%"
Set-TextValue $ws "G6" "12"

Set-TextValue $ws "D7" "1.949"
Set-TextValue $ws "E7" "-4.07%"
Set-TextValue $ws "G7" "12"

Set-TextValue $ws "E8" "-3.17%"
Set-TextValue $ws "G8" "12"

Set-TextValue $ws "D9" "0.9482"
Set-TextValue $ws "E9" "2.45%"
Set-TextValue $ws "G9" "12"

Set-TextValue $ws "D10" "0.1229"
Set-TextValue $ws "E10" "-3.25%"
Set-TextValue $ws "G10" "12"

Set-TextValue $ws "D11" "0.1975"
Set-TextValue $ws "E11" "1.07%"
Set-TextValue $ws "G11" "12"

Set-TextValue $ws "D12" "0.1008"
Set-TextValue $ws "E12" "7.14%"
Set-TextValue $ws "G12" "12"

Set-TextValue $ws "D13" "0.04436"
Set-TextValue $ws "E13" "13.39%"
Set-TextValue $ws "G13" "12"

Set-TextValue $ws "E14" "0.56%"
Set-TextValue $ws "G14" "12"

Set-TextValue $ws "D15" "0.001291"
Set-TextValue $ws "E15" "-0.94%"
Set-TextValue $ws "G15" "12"

Set-TextValue $ws "D16" "0.006074"
Set-TextValue $ws "E16" "-0.93%"
Set-TextValue $ws "G16" "12"

Set-TextValue $ws "D17" "3.473"
Set-TextValue $ws "E17" "1.02%"
Set-TextValue $ws "G17" "12"

Set-TextValue $ws "D18" "4.527"
Set-TextValue $ws "E18" "-0.11%"
Set-TextValue $ws "G18" "12"

Set-TextValue $ws "E19" "0.13%"
Set-TextValue $ws "G19" "12"

Set-TextValue $ws "D20" "8.706"
Set-TextValue $ws "E20" "4.18%"
Set-TextValue $ws "G20" "12"

Set-TextValue $ws "E21" "-0.74%"
Set-TextValue $ws "G21" "12"

Set-TextValue $ws "D22" "0.2721"
Set-TextValue $ws "E22" "2.21%"
Set-TextValue $ws "G22" "12"

Set-TextValue $ws "E23" "-0.06%"
Set-TextValue $ws "G23" "12"

Set-TextValue $ws "D24" "0.001238"
Set-TextValue $ws "E24" "-1.54%"
Set-TextValue $ws "G24" "12"

Set-TextValue $ws "D25" "0.004358"
Set-TextValue $ws "E25" "0.89%"
Set-TextValue $ws "G25" "12"

Set-TextValue $ws "D26" "0.0001262"
Set-TextValue $ws "E26" "5.04%"
Set-TextValue $ws "G26" "12"

Set-TextValue $ws "G27" "12"

Set-TextValue $ws "G28" "12"

Set-TextValue $ws "G29" "12"

Set-TextValue $ws "G30" "12"

Set-TextValue $ws "G31" "12"

Set-TextValue $ws "G32" "12"

Set-TextValue $ws "G33" "12"

Set-TextValue $ws "G34" "12"

Set-TextValue $ws "G35" "12"

Set-TextValue $ws "G36" "12"

Set-TextValue $ws "G37" "12"

Set-TextValue $ws "G38" "12"

Set-TextValue $ws "D39" "0.02809"
Set-TextValue $ws "E39" "2.14%"
Set-TextValue $ws "G39" "12"

Set-TextValue $ws "D40" "0.05899"
Set-TextValue $ws "E40" "7.00%"
Set-TextValue $ws "G40" "12"

Set-TextValue $ws "D41" "0.007928"
Set-TextValue $ws "E41" "-0.14%"
Set-TextValue $ws "G41" "12"

Set-TextValue $ws "D42" "0.1423"
Set-TextValue $ws "E42" "0.17%"
Set-TextValue $ws "G42" "12"

Set-TextValue $ws "D43" "0.009028"
Set-TextValue $ws "E43" "0.99%"
Set-TextValue $ws "G43" "12"

Set-TextValue $ws "D44" "0.002144"
Set-TextValue $ws "E44" "0.12%"
Set-TextValue $ws "G44" "12"

Set-TextValue $ws "D45" "0.009864"
Set-TextValue $ws "E45" "-16.60%"
Set-TextValue $ws "G45" "12"

Set-TextValue $ws "D46" "0.00007249"
Set-TextValue $ws "E46" "4.08%"
Set-TextValue $ws "G46" "12"

Set-TextValue $ws "D47" "0.00000000751"
Set-TextValue $ws "E47" "0.04%"
Set-TextValue $ws "G47" "12"

Set-TextValue $ws "D48" "0.003191"
Set-TextValue $ws "E48" "-0.03%"
Set-TextValue $ws "G48" "12"

Set-TextValue $ws "E49" "-0.21%"
Set-TextValue $ws "G49" "12"

Set-TextValue $ws "D50" "0.00002103"
Set-TextValue $ws "E50" "0.04%"
Set-TextValue $ws "G50" "12"

Set-TextValue $ws "D51" "0.0002002"
Set-TextValue $ws "E51" "0.04%"
Set-TextValue $ws "G51" "12"
